$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '35.393.71'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.917.32'
$ws.Range('E3').Value = '  +3.14%  '
Set-TextValue 'D5' '244.75'
$ws.Range('E5').Value = '  +2.12%  '
Set-TextValue 'D6' '0.658'
$ws.Range('E6').Value = '  +5.66%  '
$ws.Range('E7').Value = '  -0.43%  '
Set-TextValue 'D8' '41.45'
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('E9').Value = '  +7.00%  '
Set-TextValue 'D10' '52.98'
$ws.Range('E10').Value = '  +13.01%  '
Set-TextValue 'D11' '0.0716'
$ws.Range('E11').Value = '  +3.36%  '
Set-TextValue 'D12' '0.0994'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '2.194.77'
$ws.Range('E13').Value = '  +3.10%  '
$ws.Range('E14').Value = '  +5.52%  '
Set-TextValue 'D15' '0.702'
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('D16').Value = '1.901.53'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('E17').Value = '  +3.54%  '
$ws.Range('D18').Value = '35.358.88'
$ws.Range('E18').Value = '  +0.45%  '
Set-TextValue 'D19' '72.15'
$ws.Range('E19').Value = '  +3.48%  '
$ws.Range('D20').Value = '0.0₃0823'
$ws.Range('E20').Value = '  +3.54%  '
Set-TextValue 'D21' '239.80'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('E22').Value = '  +2.36%  '
Set-TextValue 'D23' '4.84'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  +1.08%  '
Set-TextValue 'D26' '2.30'
$ws.Range('E26').Value = '  +18.92%  '
Set-TextValue 'D27' '170.29'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('E28').Value = '  +5.97%  '
Set-TextValue 'D29' '18.48'
$ws.Range('E29').Value = '  +4.66%  '
$ws.Range('E30').Value = '  +2.15%  '
Set-TextValue 'D31' '4.15'
$ws.Range('E31').Value = '  +3.96%  '
Set-TextValue 'D33' '0.948'
$ws.Range('E33').Value = '  +13.61%  '
$ws.Range('E34').Value = '  -0.41%  '
Set-TextValue 'D35' '4.13'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('E36').Value = '  -4.71%  '
$ws.Range('E37').Value = '  +2.58%  '
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('E39').Value = '  +3.17%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D40' '0.0655'
$ws.Range('E40').Value = '  +14.48%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D41' '0.0209'
$ws.Range('E41').Value = '  +4.16%  '
Set-TextValue 'D42' '16.25'
$ws.Range('E42').Value = '  +9.43%  '
Set-TextValue 'D43' '90.83'
$ws.Range('E43').Value = '  +1.42%  '
$ws.Range('D44').Value = '1.344.24'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('E45').Value = '  +2.57%  '
Set-TextValue 'D46' '48.58'
$ws.Range('E46').Value = '  +39.66%  '
Set-TextValue 'D47' '2.80'
$ws.Range('E47').Value = '  +2.30%  '
$ws.Range('E48').Value = '  -0.28%  '
Set-TextValue 'D49' '6.60'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = '2.104.43'
$ws.Range('E50').Value = '  +2.99%  '
Set-TextValue 'D51' '0.0701'
$ws.Range('E51').Value = '  +3.01%  '
